$d = $word.ActiveDocument

# The five new "flowing" paragraph texts that replace the old
# "Chapter N: Brand / tagline / (Expand ...)" paragraphs (which used
# <w:br/> line breaks inside a single run).
$texts = @(
  "Toyota is known worldwide for reliability and innovation. Its Prius was a pioneer in hybrid technology, and the brand continues to lead in efficiency and quality. Toyota’s production system set global manufacturing standards.",
  "Ford revolutionized the industry with mass production. Iconic models like the Mustang represent American automotive culture. The F-series truck remains one of the best-selling vehicles worldwide.",
  "Volkswagen offers a broad lineup and global presence. Its Beetle became a cultural icon, while recent challenges like the emissions scandal highlight the complexities of modern automaking.",
  "Tesla disrupted the market by focusing exclusively on electric vehicles. Its innovations in battery technology, software, and self-driving have made it a market leader despite controversies.",
  "Ferrari and Lamborghini symbolize luxury and performance. These Italian brands combine motorsport heritage with exclusivity, appealing to enthusiasts and collectors globally."
)

# Step 1: rewrite the five existing chapter paragraphs (paragraphs 2-6,
# right after the title paragraph) with the new combined text, removing
# the <w:br/> separated runs.
for ($i = 0; $i -lt $texts.Length; $i++) {
  $p = $d.Paragraphs.Item($i + 2)
  $p.Range.Text = $texts[$i]
}

# Step 2: duplicate the block of five paragraphs six more times so the
# document ends up with seven total repeats (35 paragraphs) of the five
# brand write-ups, appended at the end of the body (before sectPr).
for ($rep = 0; $rep -lt 6; $rep++) {
  for ($i = 0; $i -lt $texts.Length; $i++) {
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = $texts[$i]
  }
}
